$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column BH: quarter "Agosto.2021" -------------------------------
# Header cell (row 1) with the same bold/centered/bordered style as the
# rest of the header row.
$ws.Range("BH1").Value = "Agosto.2021"
$ws.Range("BG1").Copy()
$ws.Range("BH1").PasteSpecial(-4122)   # xlPasteFormats

# Rows 2-73 simply carry forward the last known value (same as column BG).
$ws.Range("BG2:BG73").Copy()
$ws.Range("BH2:BH73").PasteSpecial(-4163)   # xlPasteValues

# Row 74 gets its own, newly published figure (differs from BG74).
$ws.Range("BH74").Value = 3600

# --- New row 75: quarter "01-04-2021" -----------------------------------
# The label looks like a date, so Excel would otherwise silently convert
# it to a date serial number; a leading quote forces literal text, and we
# immediately normalise the cell's format back to the plain, unstyled
# look used by the rest of column A.
$ws.Range("A75").Value = "'01-04-2021"
$ws.Range("A2").Copy()
$ws.Range("A75").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("BH75").Value = 3775

$excel.CutCopyMode = $false
